$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update total mora value (Valor Mora) and number of periods (Cant. Periodos)
$ws.Range("E11").Value = 120000
$ws.Range("F13").Value = 2

# Insert a new row for the new period (2508), pushing the signature block down
$ws.Rows("17:17").Insert()

# Duplicate the existing worker data row (16) into the freshly inserted row (17)
$ws.Range("B16:J16").Copy()
$ws.Range("B17:J17").PasteSpecial()

# New row corresponds to period 2508 (same worker, same values)
$ws.Range("E17").Value = "2508"
